# Updates cryptos list data (prices/volume %) in the "Cryptos" sheet,
# matching the GitHub Actions scheduled refresh commit.
# For cells whose new text looks like a plain number (e.g. "267.35" or
# "1.00"), we briefly force Text format so Excel keeps it as a literal
# string (preserving trailing zeros / multi-dot "thousands" price
# formatting) instead of silently converting it to a numeric value,
# then restore the "Normal" style so no visible formatting changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '43.682.80'
$ws.Cells.Item(2,5).Value = '  -0.27%  '
$ws.Cells.Item(3,4).Value = '2.315.05'
$ws.Cells.Item(3,5).Value = '  +4.52%  '
$ws.Cells.Item(4,5).Value = '  +0.22%  '
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '267.35'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = '  +1.11%  '
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '91.83'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = '  +6.12%  '
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '0.630'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = '  +2.21%  '
$ws.Cells.Item(8,5).Value = '  +0.12%  '
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.616'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = '  +1.68%  '
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '44.39'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = '  -4.14%  '
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.0933'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = '  +1.18%  '
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '7.97'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = '  +4.76%  '
$ws.Cells.Item(13,5).Value = '  +0.14%  '
$ws.Cells.Item(14,4).Value = '2.657.79'
$ws.Cells.Item(14,5).Value = '  +4.44%  '
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '15.25'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = '  +3.92%  '
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '0.853'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = '  +8.50%  '
$ws.Cells.Item(17,4).Value = '2.322.92'
$ws.Cells.Item(17,5).Value = '  +5.51%  '
$ws.Cells.Item(18,4).Value = '43.725.56'
$ws.Cells.Item(18,5).Value = '  -0.03%  '
$ws.Cells.Item(19,5).Value = '  +0.93%  '
$ws.Cells.Item(20,5).Value = '  +4.38%  '
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '71.04'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = '  +1.38%  '
$ws.Cells.Item(22,2).Value = 'BitcoinCash'
$ws.Cells.Item(22,3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '237.55'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = '  +2.33%  '
$ws.Cells.Item(23,2).Value = 'ImmutableX'
$ws.Cells.Item(23,3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '2.26'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = '  -5.24%  '
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '9.53'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = '  +7.39%  '
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '1.00'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = '  +0.11%  '
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '2.48'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = '  -1.16%  '
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '11.10'
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = '  +2.12%  '
$ws.Cells.Item(28,5).Value = '  -3.85%  '
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '2.27'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = '  +2.42%  '
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '38.70'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = '  -1.53%  '
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '22.79'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = '  +10.88%  '
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '172.63'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = '  -1.48%  '
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '0.0882'
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = '  -1.50%  '
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '5.46'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = '  +1.06%  '
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '0.125'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = '  +1.23%  '
$ws.Cells.Item(36,5).Value = '  -0.25%  '
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '4.46'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = '  +0.51%  '
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.0347'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = '  -3.00%  '
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '3.33'
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = '  +1.47%  '
$ws.Cells.Item(40,5).Value = '  +11.32%  '
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.232'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = '  +14.08%  '
$ws.Cells.Item(42,2).Value = 'ARBITRUM'
$ws.Cells.Item(42,3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '1.34'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = '  +19.47%  '
$ws.Cells.Item(43,2).Value = 'Celestia'
$ws.Cells.Item(43,3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '12.06'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = '  -2.84%  '
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '5.42'
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = '  -1.52%  '
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '60.74'
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = '  -6.80%  '
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '8.91'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = '  +6.59%  '
$ws.Cells.Item(47,5).Value = '  +2.65%  '
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '99.94'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = '  -0.19%  '
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '1.18'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = '  -0.47%  '
$ws.Cells.Item(50,4).Value = '2.542.23'
$ws.Cells.Item(50,5).Value = '  +4.68%  '
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '0.427'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = '  -3.92%  '
